$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old rows 8-10 (sender=MuSCs block superseded by restructured data)
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Rewrite rows 2-7 with the updated TPM values
# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Ephb6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.900731333333333
$ws.Cells.Item(2, 8).Value = 8.702194
$ws.Cells.Item(2, 9).Value = 0.8130494232775288
$ws.Cells.Item(2, 10).Value = 0.8130494232775289
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2226943333333333
$ws.Cells.Item(2, 14).Value = 0.668083
$ws.Cells.Item(2, 15).Value = 0.09471013227150191
$ws.Cells.Item(2, 16).Value = 0.09471013227150192
$ws.Cells.Item(2, 17).Value = 0.6459764304557778
$ws.Cells.Item(2, 18).Value = 5.813787874102
$ws.Cells.Item(2, 19).Value = 0.07700401842188309
$ws.Cells.Item(2, 20).Value = 0.07700401842188312

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Ephb6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.900731333333333
$ws.Cells.Item(3, 8).Value = 8.702194
$ws.Cells.Item(3, 9).Value = 0.8130494232775288
$ws.Cells.Item(3, 10).Value = 0.8130494232775289
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.8319233333333332
$ws.Cells.Item(3, 14).Value = 2.49577
$ws.Cells.Item(3, 15).Value = 0.3538103900551972
$ws.Cells.Item(3, 16).Value = 0.3538103900551972
$ws.Cells.Item(3, 17).Value = 2.413186079931111
$ws.Cells.Item(3, 18).Value = 21.71867471938
$ws.Cells.Item(3, 19).Value = 0.2876653335839756
$ws.Cells.Item(3, 20).Value = 0.2876653335839757

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Ephb6"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.900731333333333
$ws.Cells.Item(4, 8).Value = 8.702194
$ws.Cells.Item(4, 9).Value = 0.8130494232775288
$ws.Cells.Item(4, 10).Value = 0.8130494232775289
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.296707666666667
$ws.Cells.Item(4, 14).Value = 3.890123
$ws.Cells.Item(4, 15).Value = 0.5514794776733007
$ws.Cells.Item(4, 16).Value = 0.5514794776733009
$ws.Cells.Item(4, 17).Value = 3.761400558873556
$ws.Cells.Item(4, 18).Value = 33.852605029862
$ws.Cells.Item(4, 19).Value = 0.44838007127167
$ws.Cells.Item(4, 20).Value = 0.4483800712716702

# Row 5
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Ephb6"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.6669870000000001
$ws.Cells.Item(5, 8).Value = 2.000961
$ws.Cells.Item(5, 9).Value = 0.1869505767224711
$ws.Cells.Item(5, 10).Value = 0.1869505767224711
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2226943333333333
$ws.Cells.Item(5, 14).Value = 0.668083
$ws.Cells.Item(5, 15).Value = 0.09471013227150191
$ws.Cells.Item(5, 16).Value = 0.09471013227150192
$ws.Cells.Item(5, 17).Value = 0.148534225307
$ws.Cells.Item(5, 18).Value = 1.336808027763
$ws.Cells.Item(5, 19).Value = 0.0177061138496188
$ws.Cells.Item(5, 20).Value = 0.0177061138496188

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Ephb6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.6669870000000001
$ws.Cells.Item(6, 8).Value = 2.000961
$ws.Cells.Item(6, 9).Value = 0.1869505767224711
$ws.Cells.Item(6, 10).Value = 0.1869505767224711
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8319233333333332
$ws.Cells.Item(6, 14).Value = 2.49577
$ws.Cells.Item(6, 15).Value = 0.3538103900551972
$ws.Cells.Item(6, 16).Value = 0.3538103900551972
$ws.Cells.Item(6, 17).Value = 0.55488204833
$ws.Cells.Item(6, 18).Value = 4.99393843497
$ws.Cells.Item(6, 19).Value = 0.06614505647122156
$ws.Cells.Item(6, 20).Value = 0.06614505647122157

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Ephb6"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.6669870000000001
$ws.Cells.Item(7, 8).Value = 2.000961
$ws.Cells.Item(7, 9).Value = 0.1869505767224711
$ws.Cells.Item(7, 10).Value = 0.1869505767224711
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.296707666666667
$ws.Cells.Item(7, 14).Value = 3.890123
$ws.Cells.Item(7, 15).Value = 0.5514794776733007
$ws.Cells.Item(7, 16).Value = 0.5514794776733009
$ws.Cells.Item(7, 17).Value = 0.8648871564670001
$ws.Cells.Item(7, 18).Value = 7.783984408203001
$ws.Cells.Item(7, 19).Value = 0.1030994064016307
$ws.Cells.Item(7, 20).Value = 0.1030994064016307
